# Update "algorithm and flowchart.docx":
#   1. "if(" -> "if (" in the "If false check if(...)" paragraph, and drop the
#      w:proofErr gramStart/gramEnd markers that wrapped it.
#   2. "If false print C" -> "If false, print C." in the following paragraph,
#      and drop the w:proofErr gramStart/gramEnd markers that wrapped "false".
#
# Because w:proofErr markers aren't exposed as first-class objects on the
# Word OM, we rebuild the two affected paragraphs' XML (via Range.InsertXML)
# with the corrected run text/splits and without the proofErr elements.

$d = $word.ActiveDocument
$paras = $d.Paragraphs

$rPr32 = '<w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'

# --- Paragraph 1: "If false check if(B is greater than C ...)" ---
$targetIf = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*false check if(*") {
        $targetIf = $p
        break
    }
}

if ($targetIf -eq $null) {
    throw "Could not locate the 'if false check if(' paragraph"
}

$pPrIf = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr>' + $rPr32 + '</w:pPr>'
$runsIf = (
    ('<w:r>' + $rPr32 + '<w:t>I</w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t xml:space="preserve">f false check </w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t>if (</w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t xml:space="preserve">B is greater than C and </w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t>B is less than C) or (B is less than A or B is greater</w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t xml:space="preserve"> than C)</w:t></w:r>')
) -join ''

$xmlIf = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrIf + $runsIf + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetIf.Range.InsertXML($xmlIf)

# --- Paragraph 2: "If false print C" ---
$targetFalse = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*If *false*print C*") {
        $targetFalse = $p
        break
    }
}

if ($targetFalse -eq $null) {
    throw "Could not locate the 'If false print C' paragraph"
}

$pPrFalse = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr>' + $rPr32 + '</w:pPr>'
$runsFalse = (
    ('<w:r>' + $rPr32 + '<w:t xml:space="preserve">If </w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t>false,</w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t xml:space="preserve"> print </w:t></w:r>'),
    ('<w:r>' + $rPr32 + '<w:t>C.</w:t></w:r>')
) -join ''

$xmlFalse = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $pPrFalse + $runsFalse + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetFalse.Range.InsertXML($xmlFalse)

Write-Host "Done: if( -> if ( ; false print C -> false, print C."
